# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    2  = 191
    3  = 3012
    4  = 215
    5  = 113
    7  = 1644
    8  = 1612
    14 = 26
    19 = 14
    21 = 10
    22 = 357
    23 = 156
    24 = 95
    25 = 18
    26 = 2017
    28 = 458
    29 = 16
    30 = 174
    33 = 332
    34 = 4
    35 = 492
    36 = 6
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
